$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "24.936.73"
$ws.Range("E2").Value2 = "  +2.13%  "
$ws.Range("D3").Value2 = "1.677.48"
$ws.Range("E3").Value2 = "  +1.48%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value2 = "1.000"
$ws.Range("E4").Value2 = "  -0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "328.16"
$ws.Range("E5").Value2 = "  +7.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "0.9984"
$ws.Range("E6").Value2 = "  +0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = "0.3652"
$ws.Range("E7").Value2 = "  +0.66%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "47.16"
$ws.Range("E8").Value2 = "  -0.44%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "0.3258"
$ws.Range("E9").Value2 = "  -0.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "1.147"
$ws.Range("E10").Value2 = "  +2.66%  "
$ws.Range("E11").Value2 = "  +2.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "0.9984"
$ws.Range("E12").Value2 = "  -0.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "6.093"
$ws.Range("E13").Value2 = "  +2.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "19.73"
$ws.Range("E14").Value2 = "  +2.78%  "
$ws.Range("D15").Value2 = "1.676.95"
$ws.Range("E15").Value2 = "  +1.51%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "6.633"
$ws.Range("E16").Value2 = "  +1.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "0.00001055"
$ws.Range("E17").Value2 = "  +0.90%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "0.06604"
$ws.Range("E18").Value2 = "  +1.79%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "0.9976"
$ws.Range("E19").Value2 = "  -0.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "79.12"
$ws.Range("E20").Value2 = "  +2.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "15.96"
$ws.Range("E21").Value2 = "  +1.30%  "
$ws.Range("E22").Value2 = "  +0.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "12.52"
$ws.Range("E23").Value2 = "  +2.46%  "
$ws.Range("D24").Value2 = "24.921.62"
$ws.Range("E24").Value2 = "  +2.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "2.452"
$ws.Range("E25").Value2 = "  +0.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "2.425"
$ws.Range("E26").Value2 = "  +3.55%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "148.60"
$ws.Range("E27").Value2 = "  +1.89%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "18.81"
$ws.Range("E28").Value2 = "  +1.78%  "
$ws.Range("D29").Value2 = "1.861.78"
$ws.Range("E29").Value2 = "  +1.38%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "126.08"
$ws.Range("E30").Value2 = "  +1.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "1.190"
$ws.Range("E31").Value2 = "  +2.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "4.073"
$ws.Range("E32").Value2 = "  +0.82%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "5.780"
$ws.Range("E33").Value2 = "  +3.41%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "0.08461"
$ws.Range("E34").Value2 = "  +1.48%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "1.646"
$ws.Range("E35").Value2 = "  -1.61%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "12.27"
$ws.Range("E36").Value2 = "  +0.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "5.181"
$ws.Range("E37").Value2 = "  +0.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "0.02273"
$ws.Range("E38").Value2 = "  +2.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "0.06071"
$ws.Range("E39").Value2 = "  +0.34%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "1.231"
$ws.Range("E40").Value2 = "  +2.10%  "
$ws.Range("E41").Value2 = "  +2.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "8.295"
$ws.Range("E42").Value2 = "  +0.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "0.9975"
$ws.Range("E43").Value2 = "  -0.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "0.5976"
$ws.Range("E44").Value2 = "  +2.17%  "
$ws.Range("E45").Value2 = "  +6.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "3.839"
$ws.Range("E46").Value2 = "  +3.20%  "
$ws.Range("E47").Value2 = "  +2.80%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "125.66"
$ws.Range("E48").Value2 = "  +3.47%  "
$ws.Range("E49").Value2 = "  +1.63%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "0.07015"
$ws.Range("E50").Value2 = "  +1.75%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "1.189"
$ws.Range("E51").Value2 = "  +3.41%  "
